$wb = $excel.ActiveWorkbook
$wsAdmin = $wb.Worksheets.Item("SuperAdmin")
$ws = $wb.Worksheets.Item("Customer")

# --- Reorder column A (rows 2-5): cut A5's value and insert it above A2,
#     shifting the three values that were in A2:A4 down into A3:A5.
#     (Equivalent to Excel's "Cut A5" + "Insert Cut Cells" at A2, restricted
#     to column A - column B is untouched.)
$cut = $ws.Range("A5").Value2
$ws.Range("A5").Value2 = $ws.Range("A4").Value2
$ws.Range("A4").Value2 = $ws.Range("A3").Value2
$ws.Range("A3").Value2 = $ws.Range("A2").Value2
$ws.Range("A2").Value2 = $cut

# --- Add a hyperlink on A7 (picks up Excel's built-in "Hyperlink" cell style) ---
$ws.Hyperlinks.Add($ws.Range("A7"), "mailto:hascap814@tapi.re")

# --- Add a hyperlink on A2, but keep its original (non-hyperlink) formatting ---
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:padnunlow@addrin.uk")
$ws.Range("A6").Copy()
$ws.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

# --- Make "Customer" the active/selected sheet & tab, with a fresh selection ---
$ws.Activate()
$ws.Range("E5").Select()

# SuperAdmin keeps its own last selection untouched
$wsAdmin.Range("E6").Select()
$ws.Activate()
